# "cleaned defensive actions data"
# The sheet had a two-row header (row 1 mostly blank/merged "Unnamed: x"
# placeholders from a pandas export, row 2 the real column labels). This
# edit promotes the real labels onto row 1, removes the merges that used
# to visually fake the grouping, hides the now-redundant second header
# row (plus the blank spacer row and the summary row), and fills in the
# handful of previously-omitted zero values in the Tkl% column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the merges that used to group H1:L1 ("Tackles"), M1:P1
#     ("Challenges") and Q1:S1 under one label each, *before* writing
#     individual values into the now-freed cells (a merged range only
#     accepts a value on its top-left anchor cell). ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- Row 1: replace the pandas "Unnamed: n_level_0" / merged-group
#     placeholders with the real per-column headers. ---
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Row 2 (the old duplicate header row), row 3 (blank spacer row)
#     and row 20 (the aggregate "16 Players" summary row) are kept but
#     hidden now that row 1 carries real labels. ---
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(20).Hidden = $true

# --- Fill in the previously-missing Tkl% zeros so every player row has
#     a value in column O. ---
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("O19").Value = 0

# --- Restore the cursor to where the author left it. ---
$ws.Range("O21").Select() | Out-Null
